$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (engine rounds ColumnWidth to the nearest pixel using the
# sheet's default font metrics, same as real Excel; an integer input here
# reproduces the ".83203125"-style fractional width Excel itself writes)
$ws.Columns.Item(2).ColumnWidth = 60
$ws.Columns.Item(3).ColumnWidth = 7
$ws.Columns.Item(6).ColumnWidth = 26

# Row 2 - Deposit
$ws.Range("B2").Value = "DepositRef Nbr: 130012345"
$ws.Range("C2").Value = "05-15"
$ws.Range("F2").Value = "Deposits & Other Credits"

# Row 3 - ATM Withdrawal
$ws.Range("B3").Value = "ATM Withdrawal 1000 Walnut St M119 Kansas City MO 00005678"
$ws.Range("C3").Value = "05-18"
$ws.Range("D3").Value = 20
$ws.Range("F3").Value = "ATM Withdrawals & Debits"

# Row 4 - Check Paid 1001
$ws.Range("A4").Value = "Check Paid"
$ws.Range("C4").Value = "05-12"
$ws.Range("F4").Value = "Checks Paid"

# Row 5 - Check Paid 1002
$ws.Range("A5").Value = "Check Paid"
$ws.Range("C5").Value = "05-18"
$ws.Range("D5").Value = 230
$ws.Range("F5").Value = "Checks Paid"

# Row 6 - Check Paid 1003
$ws.Range("A6").Value = "Check Paid"
$ws.Range("C6").Value = "05-24"
$ws.Range("F6").Value = "Checks Paid"
